# edit.ps1 - applies the V103 "Auswertung" content expansion described in the diff.
# Strategy: use Range.InsertXML with self-contained WordprocessingML fragments
# (the same approach Word's COM InsertXML accepts) to (1) replace paragraph 1
# ("Auswertung") with a bold+underlined version, (2) replace paragraph 2
# ("Runder Stab ...") with the updated, underlined heading text, and
# (3) insert the large block of new report paragraphs right before the
# trailing bookmark paragraph, whose own trailing run+bookmark Word will
# merge into automatically.

$d = $word.ActiveDocument

# --- 1. "Auswertung" heading becomes bold (keeps its existing underline) ---
$para1Xml = @'
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>Auswertung</w:t></w:r></w:p></w:body></w:document>
'@

$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertXML($para1Xml)

# --- 2. "Runder Stab, einseitige Einspannung" heading gets updated text ---
#        ("(golden)" inserted) and becomes fully underlined, still numbered ---
$para2Xml = @'
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Runder Stab (golden),</w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> einseitige Einspannung</w:t></w:r></w:p></w:body></w:document>
'@

$p2 = $d.Paragraphs.Item(2)
$p2.Range.InsertXML($para2Xml)

# --- 3. Insert all of the new "Auswertung" body content before the final
#        (bookmarked) paragraph. The fragment's last paragraph mark merges
#        into the existing trailing paragraph, carrying the _GoBack bookmark
#        forward exactly like the target document. ---
$middleXml = @'
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>    <w:p>
      <w:r>
        <w:t xml:space="preserve">Die eingespannte Länge L des Stabes und die Masse m des an das Ende des Stabes angehängte Gewicht betragen </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>… .</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Der Stab hat einen Durchmesser von … .</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Bei der linearen Ausgleichsrechnung wird D(x) auf der y-Achse und Lx^3-x^3/3 auf der x-Achse aufgetragen. D(x) ist dabei gemäß D(x)=D_m(x)-D_0(x) zu berechnen.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> NUMMERIERT</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Die Werte sind Tabelle 1 zu entnehmen.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Tabelle 1</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Die von Python bestimmten Werte für a und b sind dann:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>…</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">… </w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Mit y(x)=a*x+b ist der Elastizitätsmodul durch </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">E =… </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Dabei gilt für I </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">wegen einer runden Querschnittsfläche </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">nach </w:t>
      </w:r>
      <w:r>
        <w:t>(</w:t>
      </w:r>
      <w:r>
        <w:t>QUELLE</w:t>
      </w:r>
      <w:r>
        <w:t>)</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> I = … , für g wird 9.81 m/s^2 verwendet.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Daraus folgt für E dann </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>…</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Der verwendete Stab war golden, sodass eventuell von Messing ausgegangen werden kann. </w:t>
      </w:r>
      <w:r>
        <w:t>Der Literaturwert ist … (QUELLE2), sodass sich eine Abweichung von … ergibt.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Listenabsatz"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:rPr>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>Eckiger Stab (grau), einseitige Einspannung</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Für die eingespannte Länge L und die Masse des Gewichts gilt hier</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>… .</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Die Seitenlänge wurde bestimmt zu</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> … .</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Wieder</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> wird D(x) auf der y-Achse und Lx^3-x^3/3 auf der x-Achse aufgetragen. </w:t>
      </w:r>
      <w:r>
        <w:t>D(x) berechnet sich erneut gemäß (NUMMERIERT).</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Die Werte sind Tabelle 2</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> zu entnehmen.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Tabelle 2</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Die Parameter a und b wurden von Python berechnet:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>…</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">… </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t xml:space="preserve">Erneut gilt für der E-Modul </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">E =… </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Es liegt eine quadratische Querschnittsfläche vor, sodass für I nach </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">QUELLE I = … </w:t>
      </w:r>
      <w:r>
        <w:t>gilt</w:t>
      </w:r>
      <w:r>
        <w:t>, für g wird</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> erneut</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> 9.81 m/s^2 verwendet.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Dann ist das E-Modul</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>…</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Der Stab hat eine graue Farbe, daher wird als Material Eisen geschätzt. Für Eisen ist der Literaturwert … (QUELLE2). Daher ist die Abweichung der Messung … .</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Listenabsatz"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:rPr>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>Eckiger Stab (grau), beidseitige Einspannung</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">L und m betragen bei der beidseitigen Einspannung </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>…</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Die Seitenlänge ist erneut … .</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Gemäß Versuchsanleitung wird D(x), bestimmt gemäß (NUMMERIERT), gegen 3L^2x-4x</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>^3 aufgetragen. Es sei angemerkt, dass deswegen die Hälfte der Werte, nämlich genau die bis zum angehängten Gewicht, verloren gehen, da für die zweite Hälfte die Gleichung (NUMMERIERT AUS THEORIETEIL) gilt.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Die aufgetragenen Werte sind der untenstehenden Tabelle 3 zu entnehmen</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Tabelle 3</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Für a und b gilt nach linearer Regression durch Python</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>…</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>…</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">In dem Fall der beidseitigen Einspannung gilt für E </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>…,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>wobei für I erneut … verwendet wird.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Nach Einsetzen erhält man für E </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>… .</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Die Abweichung beträgt bei dem Literaturwert Elit = … …%.</w:t>
      </w:r>
</w:p></w:body></w:document>
'@

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertionPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$insertionPoint.InsertXML($middleXml)

Write-Output "Applied Auswertung edits; paragraph count is now $($d.Paragraphs.Count)"
